{"js": "// \"add sentence on credit card fees\"\n//\n// 1. Insert a new paragraph - \"If you use a credit card to pay, you may\n//    also be charged a nonrefundable 3% processing fee.\" - right before\n//    the paragraph that starts \"If you can't afford to pay the filing\n//    fee...\", reusing that paragraph's formatting (same pPr/rPr).\n// 2. Tidy up the \"Collect\" + \"i\" + \"ng Your Judgment\" hyperlink, which was\n//    split across three runs, into a single \"Collecting Your Judgment\"\n//    run (keeping the hyperlink's blue/underlined look).\nconst body = context.document.body;\n\n// --- 1. New \"credit card\" sentence -----------------------------------\nconst feeWaiverResults = body.search(\n  \"If you can\\u2019t afford to pay the filing fee\",\n  { matchCase: true }\n);\nfeeWaiverResults.load(\"items\");\nawait context.sync();\n\nif (feeWaiverResults.items.length > 0) {\n  const feeWaiverRange = feeWaiverResults.items[0];\n  const feeWaiverParagraphs = feeWaiverRange.paragraphs;\n  feeWaiverParagraphs.load(\"items\");\n  await context.sync();\n\n  const feeWaiverParagraph = feeWaiverParagraphs.items[0];\n  feeWaiverParagraph.insertParagraph(\n    \"If you use a credit card to pay, you may also be charged a nonrefundable 3% processing fee.\",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n}\n\n// --- 2. Merge the \"Collecting Your Judgment\" link text into one run ---\nconst judgmentResults = body.search(\"Collecting Your Judgment\", {\n  matchCase: true,\n});\njudgmentResults.load(\"items\");\nawait context.sync();\n\nif (judgmentResults.items.length > 0) {\n  const judgmentRange = judgmentResults.items[0];\n  judgmentRange.font.load(\"color,underline\");\n  await context.sync();\n\n  const savedColor = judgmentRange.font.color;\n  const savedUnderline = judgmentRange.font.underline;\n\n  // Replacing with identical text is a no-op for the engine (the three\n  // runs would be left untouched), so swap in a placeholder first to\n  // force a real rewrite, then restore the real text and formatting.\n  judgmentRange.insertText(\"ZZZPLACEHOLDERZZZ\", Word.InsertLocation.replace);\n  await context.sync();\n\n  const placeholderResults = context.document.body.search(\n    \"ZZZPLACEHOLDERZZZ\",\n    { matchCase: true }\n  );\n  placeholderResults.load(\"items\");\n  await context.sync();\n\n  if (placeholderResults.items.length > 0) {\n    const placeholderRange = placeholderResults.items[0];\n    placeholderRange.insertText(\n      \"Collecting Your Judgment\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n\n    const finalResults = context.document.body.search(\n      \"Collecting Your Judgment\",\n      { matchCase: true }\n    );\n    finalResults.load(\"items\");\n    await context.sync();\n\n    if (finalResults.items.length > 0) {\n      const finalRange = finalResults.items[0];\n      finalRange.font.color = savedColor;\n      finalRange.font.underline = savedUnderline;\n      await context.sync();\n    }\n  }\n}\n", "ps1": "# \"add sentence on credit card fees\"\n#\n# 1. Insert a new paragraph - \"If you use a credit card to pay, you may\n#    also be charged a nonrefundable 3% processing fee.\" - right before\n#    the paragraph that starts \"If you can't afford to pay the filing\n#    fee...\", reusing that paragraph's formatting.\n# 2. Tidy up the \"Collect\" + \"i\" + \"ng Your Judgment\" hyperlink, which was\n#    split across three runs, into a single \"Collecting Your Judgment\"\n#    run (keeping the hyperlink's blue/underlined look).\n\n$d = $word.ActiveDocument\n\n# --- 1. New \"credit card\" sentence ------------------------------------\n$feeRange = $d.Content\n$feeRange.Find.ClearFormatting()\n$feeFound = $feeRange.Find.Execute(\"If you can\" + [char]0x2019 + \"t afford to pay the filing fee\")\nif ($feeFound) {\n    $feeRange.Collapse(1)  # wdCollapseStart\n    $feeRange.InsertBefore(\"If you use a credit card to pay, you may also be charged a nonrefundable 3% processing fee.\" + [char]13)\n}\n\n# --- 2. Merge the \"Collecting Your Judgment\" link text into one run ---\n$judgRange = $d.Content\n$judgRange.Find.ClearFormatting()\n$judgFound = $judgRange.Find.Execute(\"Collecting Your Judgment\")\nif ($judgFound) {\n    # Replacing with identical text is a no-op for the engine (the three\n    # runs would be left untouched), so swap in a placeholder first to\n    # force a real rewrite, then restore the real text and formatting.\n    $judgRange.Text = \"ZZZPLACEHOLDERZZZ\"\n\n    $phRange = $d.Content\n    $phRange.Find.ClearFormatting()\n    $phFound = $phRange.Find.Execute(\"ZZZPLACEHOLDERZZZ\")\n    if ($phFound) {\n        $phRange.Text = \"Collecting Your Judgment\"\n        $phRange.Font.Color = 12673797  # RGB(0x05, 0x63, 0xC1) packed as BGR\n        $phRange.Font.Underline = 1     # wdUnderlineSingle\n    }\n}\n\nWrite-Output \"done\"\n"}
